# 10-03-2018 : BugFixed->Employee CostCode->Add New
# Adds a new "BASIC" cost-code column (AF) and a new data row (row 3)
# to the EmpMaster_Upload template's Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column AF: header (styled like the other yellow header cells) -
$ws.Range("AF1").Value = "BASIC"
$ws.Range("AF1").Interior.Color = 65535
$ws.Range("AF1").HorizontalAlignment = -4108
$ws.Range("AF2").Value = 0

# --- New row 3: mirrors row 2, with a couple of values changed ---------
# Numeric-looking text is entered with a leading apostrophe so Excel
# stores it as text (re-using the existing shared strings); the
# quote-prefix formatting that trick applies is stripped right after
# (per cell, so empty neighbour cells are never touched) so the cells
# end up with plain (unstyled) text instead of inheriting row 2's
# quote-prefixed display style.
$ws.Range("A3").Value = 9999999
$ws.Range("B3").Value = "CONT"
$ws.Range("C3").Value = "'001"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = "TEST EMP"
$ws.Range("E3").Value = "TEST FATHER"
$ws.Range("F3").Value = "M"
$ws.Range("G3").Value = "Y"
$ws.Range("H3").Value = "'1979-06-19"
$ws.Range("H3").ClearFormats()
$ws.Range("I3").Value = "'2017-11-10"
$ws.Range("I3").ClearFormats()
$ws.Range("J3").Value = "SUN"
$ws.Range("K3").Value = "Y"
$ws.Range("L3").Value = "N"
$ws.Range("M3").Value = "Y"
$ws.Range("N3").Value = "Y"
$ws.Range("O3").Value = "GN"
$ws.Range("P3").Value = "'009999"
$ws.Range("P3").ClearFormats()
$ws.Range("Q3").Value = "NRG"
$ws.Range("R3").Value = "'01"
$ws.Range("R3").ClearFormats()
$ws.Range("S3").Value = "'001"
$ws.Range("S3").ClearFormats()
$ws.Range("T3").Value = "'001"
$ws.Range("T3").ClearFormats()
$ws.Range("U3").Value = "'001"
$ws.Range("U3").ClearFormats()
$ws.Range("V3").Value = "'012"
$ws.Range("V3").ClearFormats()
$ws.Range("W3").Value = "'012"
$ws.Range("W3").ClearFormats()
$ws.Range("X3").Value = "'001"
$ws.Range("X3").ClearFormats()
$ws.Range("Y3").Value = "'001"
$ws.Range("Y3").ClearFormats()
$ws.Range("AA3").Value = "'99999"
$ws.Range("AA3").ClearFormats()
$ws.Range("AB3").Value = "PACEHCEM01"

# AC3 sits in a column that carries a column-level style, so it is reset
# via a named style instead of ClearFormats, then centered (s="1" in the
# source edit - no quote-prefix, just centered alignment).
$ws.Range("AC3").Value = "'123123123123"
$ws.Range("AC3").Style = "Normal"
$ws.Range("AC3").HorizontalAlignment = -4108

$ws.Range("AF3").Value = 320

# --- Selection ends on the newly added row ------------------------------
$ws.Range("A3").Select()
